$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.360.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.879.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.25%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7195"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.40%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07994"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.69%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3142"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.13%  "

$ws.Range("E10").Value = "  -0.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08152"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.874.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.63%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "94.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.229"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7101"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.98%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.400"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008455"
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.369.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "248.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.14%  "

$ws.Range("E20").Value = "  +0.45%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.124.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.743"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.16%  "

$ws.Range("E24").Value = "  +0.23%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1597"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.42%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.066"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.38%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.97%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.502"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.414"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.16%  "

$ws.Range("E31").Value = "  -0.83%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.220"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05339"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.75%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.938"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7562"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.179"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.702"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.61%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01887"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.266.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.763"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.97%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.438"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.58%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "113.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.36%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9061"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.60%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "74.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.89%  "

$ws.Range("E45").Value = "  +0.15%  "

$ws.Range("E46").Value = "  +4.69%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.023.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.47%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.801"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.12%  "

$ws.Range("E49").Value = "  -0.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.488"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.40%  "

$ws.Range("E51").Value = "  +0.15%  "

